{"js": "// Merge the split \"<id>...</id>\" runs back into single runs with\n// updated/normalized text, as described by the diff:\n//   <id>p126v_3</id>   (3 runs)      -> <id>p126v_3</id>   (1 run)\n//   <id>p127r_a1</id>  (3 runs)      -> <id>p127r_1</id>   (1 run)\n// search() matches across run boundaries on the plain-text content, and\n// insertText(..., Replace) collapses the matched range into a single run\n// that inherits the formatting of the first run in the match - which is\n// exactly the Courier New / 7f6000 / 18pt \"tag\" styling used for the\n// <id> and </id> runs in this document.\n\nconst body = context.document.body;\n\nconst firstIdResults = body.search(\"<id>p126v_3</id>\", { matchCase: true });\nfirstIdResults.load(\"items\");\nawait context.sync();\n\nfor (const r of firstIdResults.items) {\n  r.insertText(\"<id>p126v_3</id>\", Word.InsertLocation.replace);\n}\n\nconst secondIdResults = body.search(\"<id>p127r_a1</id>\", { matchCase: true });\nsecondIdResults.load(\"items\");\nawait context.sync();\n\nfor (const r of secondIdResults.items) {\n  r.insertText(\"<id>p127r_1</id>\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Merge the split \"<id>...</id>\" runs back into single runs with\n# updated/normalized text, as described by the commit diff:\n#   <id>p126v_3</id>   (3 runs)  -> <id>p126v_3</id>   (1 run)\n#   <id>p127r_a1</id>  (3 runs)  -> <id>p127r_1</id>   (1 run)\n#\n# Word's Find/Replace matches across run boundaries on the plain-text\n# content, and replacing the whole matched range collapses it into a\n# single run that inherits the formatting of the first run in the match\n# - i.e. the Courier New / color 7f6000 / 18pt \"tag\" styling already used\n# for the surrounding <id> / </id> runs in this document.\n\n$d = $word.ActiveDocument\n\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Execute(\n    \"<id>p126v_3</id>\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"<id>p126v_3</id>\",\n    2\n)\n\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Execute(\n    \"<id>p127r_a1</id>\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"<id>p127r_1</id>\",\n    2\n)\n"}
